$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 9560
$ws.Range("E3").Value = 15930
$ws.Range("E4").Value = 7689
$ws.Range("E5").Value = 4755
$ws.Range("E6").Value = 17670
$ws.Range("E7").Value = 13068
$ws.Range("E8").Value = 12517
$ws.Range("E9").Value = 3154
$ws.Range("E10").Value = 13002
$ws.Range("E11").Value = 12217
$ws.Range("E12").Value = 10818
$ws.Range("E13").Value = 1894
